$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2297.1667  # H19: 2356.6 -> 2297.1667
$ws.Cells.Item(19, 9).Value = 2594.3333  # I19: 2891.5 -> 2594.3333
$ws.Cells.Item(19, 11).Value = 2594.3333  # K19: 2891.5 -> 2594.3333
$ws.Cells.Item(19, 13).Value = -2419.3333  # M19: -2716.5 -> -2419.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 911290.8  # H43: 835508.2 -> 911290.8
$ws.Cells.Item(43, 9).Value = 2257  # I43: 2212.25 -> 2257
$ws.Cells.Item(43, 11).Value = 2257  # K43: 2212.25 -> 2257
$ws.Cells.Item(43, 13).Value = -2188  # M43: -2143.25 -> -2188

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 116.94444  # H55: 109.95 -> 116.94444
$ws.Cells.Item(55, 9).Value = 172.5  # I55: 158.5 -> 172.5
$ws.Cells.Item(55, 10).Value = 101.07143  # J55: 97.8125 -> 101.07143
$ws.Cells.Item(55, 11).Value = 172.5  # K55: 158.5 -> 172.5
$ws.Cells.Item(55, 12).Value = 101.07143  # L55: 97.8125 -> 101.07143
$ws.Cells.Item(55, 13).Value = 41.5  # M55: 55.5 -> 41.5
$ws.Cells.Item(55, 14).Value = -529.07143  # N55: -525.8125 -> -529.07143

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 1549.3334  # H86: 1680.2 -> 1549.3334
$ws.Cells.Item(86, 9).Value = 898  # I86: 901 -> 898
$ws.Cells.Item(86, 11).Value = 898  # K86: 901 -> 898
$ws.Cells.Item(86, 13).Value = 225  # M86: 222 -> 225

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 1549.3334  # H89: 1680.2 -> 1549.3334
$ws.Cells.Item(89, 9).Value = 898  # I89: 901 -> 898
$ws.Cells.Item(89, 11).Value = 4490  # K89: 4505 -> 4490
$ws.Cells.Item(89, 13).Value = 1126  # M89: 1111 -> 1126

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 5123.125  # H113: 5160 -> 5123.125
$ws.Cells.Item(113, 9).Value = 5331.5  # I113: 5380.6665 -> 5331.5
$ws.Cells.Item(113, 11).Value = 5331.5  # K113: 5380.6665 -> 5331.5
$ws.Cells.Item(113, 13).Value = -2077.5  # M113: -2126.6665 -> -2077.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 7772.778  # H116: 7776.222 -> 7772.778
$ws.Cells.Item(116, 9).Value = 7355.5  # I116: 7359.375 -> 7355.5
$ws.Cells.Item(116, 11).Value = 7355.5  # K116: 7359.375 -> 7355.5
$ws.Cells.Item(116, 13).Value = -3913.5  # M116: -3917.375 -> -3913.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 11236.024  # H137: 11334.775 -> 11236.024
$ws.Cells.Item(137, 9).Value = 11773.3545  # I137: 11926.767 -> 11773.3545
$ws.Cells.Item(137, 10).Value = 9570.299999999999  # J137: 9558.799999999999 -> 9570.299999999999
$ws.Cells.Item(137, 11).Value = 35320.0635  # K137: 35780.301 -> 35320.0635
$ws.Cells.Item(137, 12).Value = 28710.9  # L137: 28676.4 -> 28710.9
$ws.Cells.Item(137, 13).Value = -32770.0635  # M137: -33230.301 -> -32770.0635
$ws.Cells.Item(137, 14).Value = -33810.89999999999  # N137: -33776.39999999999 -> -33810.89999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 5146.5293  # H138: 5362.1562 -> 5146.5293
$ws.Cells.Item(138, 9).Value = 1217  # I138: 1234.6111 -> 1217
$ws.Cells.Item(138, 10).Value = 10123.934  # J138: 10669 -> 10123.934
$ws.Cells.Item(138, 11).Value = 3651  # K138: 3703.8333 -> 3651
$ws.Cells.Item(138, 12).Value = 30371.802  # L138: 32007 -> 30371.802
$ws.Cells.Item(138, 13).Value = 1489  # M138: 1436.1667 -> 1489
$ws.Cells.Item(138, 14).Value = -40651.802  # N138: -42287 -> -40651.802

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 549.7273  # H5: 659.44446 -> 549.7273
$ws.Cells.Item(5, 9).Value = 605.625  # I5: 788.8333 -> 605.625
$ws.Cells.Item(5, 11).Value = 605.625  # K5: 788.8333 -> 605.625
$ws.Cells.Item(5, 13).Value = -493.625  # M5: -676.8333 -> -493.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5012.9487  # H32: 4822.073 -> 5012.9487
$ws.Cells.Item(32, 9).Value = 4684.3423  # I32: 4505.125 -> 4684.3423
$ws.Cells.Item(32, 11).Value = 4684.3423  # K32: 4505.125 -> 4684.3423
$ws.Cells.Item(32, 13).Value = -4397.3423  # M32: -4218.125 -> -4397.3423

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(51, 8).Value = 0  # H51: 38495 -> 0
$ws.Cells.Item(51, 10).Value = 0  # J51: 38495 -> 0
$ws.Cells.Item(51, 12).Value = 0  # L51: 38495 -> 0
$ws.Cells.Item(51, 14).ClearContents()  # N51: remove (was -40007)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3947.5  # H61: 3935.625 -> 3947.5
$ws.Cells.Item(61, 9).Value = 3497.5  # I61: 3000 -> 3497.5
$ws.Cells.Item(61, 11).Value = 3497.5  # K61: 3000 -> 3497.5
$ws.Cells.Item(61, 13).Value = -3285.5  # M61: -2788 -> -3285.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 3723.8572  # H97: 3974.1667 -> 3723.8572
$ws.Cells.Item(97, 9).Value = 3659.6  # I97: 4019 -> 3659.6
$ws.Cells.Item(97, 11).Value = 3659.6  # K97: 4019 -> 3659.6
$ws.Cells.Item(97, 13).Value = -3163.6  # M97: -3523 -> -3163.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 44588.4  # H132: 47245.03 -> 44588.4
$ws.Cells.Item(132, 9).Value = 4935.5713  # I132: 5257.231 -> 4935.5713
$ws.Cells.Item(132, 11).Value = 14806.7139  # K132: 15771.693 -> 14806.7139
$ws.Cells.Item(132, 13).Value = -12276.7139  # M132: -13241.693 -> -12276.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3947.5  # H136: 3935.625 -> 3947.5
$ws.Cells.Item(136, 9).Value = 3497.5  # I136: 3000 -> 3497.5
$ws.Cells.Item(136, 11).Value = 10492.5  # K136: 9000 -> 10492.5
$ws.Cells.Item(136, 13).Value = -7942.5  # M136: -6450 -> -7942.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 549.7273  # H4: 659.44446 -> 549.7273
$ws.Cells.Item(4, 9).Value = 605.625  # I4: 788.8333 -> 605.625
$ws.Cells.Item(4, 11).Value = 605.625  # K4: 788.8333 -> 605.625
$ws.Cells.Item(4, 13).Value = -490.625  # M4: -673.8333 -> -490.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 58944.445  # H86: 78062.5 -> 58944.445
$ws.Cells.Item(86, 9).Value = 39375  # I86: 63625 -> 39375
$ws.Cells.Item(86, 10).Value = 74600  # J86: 92500 -> 74600
$ws.Cells.Item(86, 11).Value = 39375  # K86: 63625 -> 39375
$ws.Cells.Item(86, 12).Value = 74600  # L86: 92500 -> 74600
$ws.Cells.Item(86, 13).Value = -38252  # M86: -62502 -> -38252
$ws.Cells.Item(86, 14).Value = -76846  # N86: -94746 -> -76846

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 58944.445  # H89: 78062.5 -> 58944.445
$ws.Cells.Item(89, 9).Value = 39375  # I89: 63625 -> 39375
$ws.Cells.Item(89, 10).Value = 74600  # J89: 92500 -> 74600
$ws.Cells.Item(89, 11).Value = 196875  # K89: 318125 -> 196875
$ws.Cells.Item(89, 12).Value = 373000  # L89: 462500 -> 373000
$ws.Cells.Item(89, 13).Value = -191259  # M89: -312509 -> -191259
$ws.Cells.Item(89, 14).Value = -384232  # N89: -473732 -> -384232

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1950.7  # H94: 2027.0526 -> 1950.7
$ws.Cells.Item(94, 9).Value = 1022.44446  # I94: 1087.75 -> 1022.44446
$ws.Cells.Item(94, 11).Value = 1022.44446  # K94: 1087.75 -> 1022.44446
$ws.Cells.Item(94, 13).Value = -571.44446  # M94: -636.75 -> -571.44446

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1718.4286  # H134: 1383.3529 -> 1718.4286
$ws.Cells.Item(134, 9).Value = 1642.1818  # I134: 1301.1333 -> 1642.1818
$ws.Cells.Item(134, 10).Value = 1998  # J134: 2000 -> 1998
$ws.Cells.Item(134, 11).Value = 4926.5454  # K134: 3903.3999 -> 4926.5454
$ws.Cells.Item(134, 12).Value = 5994  # L134: 6000 -> 5994
$ws.Cells.Item(134, 13).Value = -2391.5454  # M134: -1368.3999 -> -2391.5454
$ws.Cells.Item(134, 14).Value = -11064  # N134: -11070 -> -11064

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 86999  # H52: 86962 -> 86999
$ws.Cells.Item(52, 10).Value = 86999  # J52: 86962 -> 86999
$ws.Cells.Item(52, 12).Value = 86999  # L52: 86962 -> 86999
$ws.Cells.Item(52, 14).Value = -87587  # N52: -87550 -> -87587

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 30560212  # H58: 34379988 -> 30560212
$ws.Cells.Item(58, 10).Value = 91668000  # J58: 137501000 -> 91668000
$ws.Cells.Item(58, 12).Value = 91668000  # L58: 137501000 -> 91668000
$ws.Cells.Item(58, 14).Value = -91668406  # N58: -137501406 -> -91668406

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3601.5417  # H134: 3774.4092 -> 3601.5417
$ws.Cells.Item(134, 9).Value = 2918.2632  # I134: 3061.5881 -> 2918.2632
$ws.Cells.Item(134, 11).Value = 8754.7896  # K134: 9184.764299999999 -> 8754.7896
$ws.Cells.Item(134, 13).Value = -6219.7896  # M134: -6649.764299999999 -> -6219.7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 30560212  # H136: 34379988 -> 30560212
$ws.Cells.Item(136, 10).Value = 91668000  # J136: 137501000 -> 91668000
$ws.Cells.Item(136, 12).Value = 275004000  # L136: 412503000 -> 275004000
$ws.Cells.Item(136, 14).Value = -275009100  # N136: -412508100 -> -275009100

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(24, 8).Value = 164.5  # H24: 4 -> 164.5
$ws.Cells.Item(24, 9).Value = 152.66667  # I24: 0 -> 152.66667
$ws.Cells.Item(24, 10).Value = 200  # J24: 4 -> 200
$ws.Cells.Item(24, 11).Value = 458.00001  # K24: 0 -> 458.00001
$ws.Cells.Item(24, 12).Value = 600  # L24: 12 -> 600
$ws.Cells.Item(24, 13).Value = -228.00001  # M24: None -> -228.00001
$ws.Cells.Item(24, 14).Value = -1060  # N24: -472 -> -1060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 2497  # H22: 10000 -> 2497
$ws.Cells.Item(22, 9).Value = 0  # I22: 10000 -> 0
$ws.Cells.Item(22, 10).Value = 2497  # J22: 0 -> 2497
$ws.Cells.Item(22, 11).Value = 0  # K22: 10000 -> 0
$ws.Cells.Item(22, 12).Value = 2497  # L22: 0 -> 2497
$ws.Cells.Item(22, 13).ClearContents()  # M22: remove (was -9471)
$ws.Cells.Item(22, 14).Value = -3555  # N22: None -> -3555

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7768.923  # H70: 8008.0835 -> 7768.923
$ws.Cells.Item(70, 9).Value = 4970.5713  # I70: 4982.5 -> 4970.5713
$ws.Cells.Item(70, 11).Value = 4970.5713  # K70: 4982.5 -> 4970.5713
$ws.Cells.Item(70, 13).Value = -4700.5713  # M70: -4712.5 -> -4700.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 7768.923  # H73: 8008.0835 -> 7768.923
$ws.Cells.Item(73, 9).Value = 4970.5713  # I73: 4982.5 -> 4970.5713
$ws.Cells.Item(73, 11).Value = 4970.5713  # K73: 4982.5 -> 4970.5713
$ws.Cells.Item(73, 13).Value = -4034.5713  # M73: -4046.5 -> -4034.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4123.5  # H80: 4709.1055 -> 4123.5
$ws.Cells.Item(80, 9).Value = 2890.9092  # I80: 3225.6667 -> 2890.9092
$ws.Cells.Item(80, 10).Value = 5356.091  # J80: 6044.2 -> 5356.091
$ws.Cells.Item(80, 11).Value = 2890.9092  # K80: 3225.6667 -> 2890.9092
$ws.Cells.Item(80, 12).Value = 5356.091  # L80: 6044.2 -> 5356.091
$ws.Cells.Item(80, 13).Value = -1892.9092  # M80: -2227.6667 -> -1892.9092
$ws.Cells.Item(80, 14).Value = -7352.091  # N80: -8040.2 -> -7352.091

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 4123.5  # H83: 4709.1055 -> 4123.5
$ws.Cells.Item(83, 9).Value = 2890.9092  # I83: 3225.6667 -> 2890.9092
$ws.Cells.Item(83, 10).Value = 5356.091  # J83: 6044.2 -> 5356.091
$ws.Cells.Item(83, 11).Value = 14454.546  # K83: 16128.3335 -> 14454.546
$ws.Cells.Item(83, 12).Value = 26780.455  # L83: 30221 -> 26780.455
$ws.Cells.Item(83, 13).Value = -9462.546  # M83: -11136.3335 -> -9462.546
$ws.Cells.Item(83, 14).Value = -36764.455  # N83: -40205 -> -36764.455

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 0  # H97: 1716.6923 -> 0
$ws.Cells.Item(97, 9).Value = 0  # I97: 756.4286 -> 0
$ws.Cells.Item(97, 10).Value = 0  # J97: 2837 -> 0
$ws.Cells.Item(97, 11).Value = 0  # K97: 756.4286 -> 0
$ws.Cells.Item(97, 12).Value = 0  # L97: 2837 -> 0
$ws.Cells.Item(97, 13).ClearContents()  # M97: remove (was -260.4286)
$ws.Cells.Item(97, 14).ClearContents()  # N97: remove (was -3829)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2639.3  # H126: 2741.8572 -> 2639.3
$ws.Cells.Item(126, 9).Value = 2656.2856  # I126: 2848.5 -> 2656.2856
$ws.Cells.Item(126, 11).Value = 7968.8568  # K126: 8545.5 -> 7968.8568
$ws.Cells.Item(126, 13).Value = -5498.8568  # M126: -6075.5 -> -5498.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2073.5  # H132: 3373.6667 -> 2073.5
$ws.Cells.Item(132, 9).Value = 2217  # I132: 2560.6667 -> 2217
$ws.Cells.Item(132, 10).Value = 1499.5  # J132: 4999.6665 -> 1499.5
$ws.Cells.Item(132, 11).Value = 6651  # K132: 7682.000100000001 -> 6651
$ws.Cells.Item(132, 12).Value = 4498.5  # L132: 14998.9995 -> 4498.5
$ws.Cells.Item(132, 13).Value = -4121  # M132: -5152.000100000001 -> -4121
$ws.Cells.Item(132, 14).Value = -9558.5  # N132: -20058.9995 -> -9558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1971.7106  # H132: 1990.0541 -> 1971.7106
$ws.Cells.Item(132, 9).Value = 1777.6  # I132: 1797.7916 -> 1777.6
$ws.Cells.Item(132, 11).Value = 5332.799999999999  # K132: 5393.3748 -> 5332.799999999999
$ws.Cells.Item(132, 13).Value = -2802.799999999999  # M132: -2863.3748 -> -2802.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 7698.8  # H136: 5345.2856 -> 7698.8
$ws.Cells.Item(136, 9).Value = 4248.5  # I136: 2930.1177 -> 4248.5
$ws.Cells.Item(136, 10).Value = 21500  # J136: 15609.75 -> 21500
$ws.Cells.Item(136, 11).Value = 12745.5  # K136: 8790.3531 -> 12745.5
$ws.Cells.Item(136, 12).Value = 64500  # L136: 46829.25 -> 64500
$ws.Cells.Item(136, 13).Value = -10195.5  # M136: -6240.3531 -> -10195.5
$ws.Cells.Item(136, 14).Value = -69600  # N136: -51929.25 -> -69600

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 7249128.5  # H126: 7578589 -> 7249128.5
$ws.Cells.Item(126, 9).Value = 8774900  # I126: 9262339 -> 8774900
$ws.Cells.Item(126, 11).Value = 26324700  # K126: 27787017 -> 26324700
$ws.Cells.Item(126, 13).Value = -26322230  # M126: -27784547 -> -26322230
